# Generate Report for Handback
# Updates the "zh-cn" and "de-de" localization-status sheets: row 8
# (af351dea-ea3d-4126-985b-09a31b54063f) now has a generated handback
# report -> populate "Latest Target File", "Latest Handback DateTime",
# and "Error Detail" (the handback version is stale), add a hyperlink on
# the new "Latest Target File" cell, and widen columns I (Latest Target
# File) and P (Error Detail) to fit the new content.

$wb = $excel.ActiveWorkbook

$warning = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9124a58a43edc3cf75f29e1d67c99b86a7adcf94/e2e/af351dea-ea3d-4126-985b-09a31b54063f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96587e4c037f0f7ea6a1c6828adbe8bb60db6e7b/e2e/af351dea-ea3d-4126-985b-09a31b54063f.md."

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$LocaleRepo,
        [string]$TargetXlf,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File (I8) becomes a hyperlink to the handback markdown,
    # same way A8 / I2 already link out to their respective ".md" files.
    $ws.Hyperlinks.Add(
        $ws.Range("I8"),
        "https://github.com/OpenLocalizationTestOrg/$LocaleRepo/blob/96587e4c037f0f7ea6a1c6828adbe8bb60db6e7b/e2e/af351dea-ea3d-4126-985b-09a31b54063f.md",
        "",
        "",
        "af351dea-ea3d-4126-985b-09a31b54063f.md"
    ) | Out-Null

    # Latest Handback File
    $ws.Range("J8").Value = $TargetXlf

    # Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # Error Detail - the handback was generated from a stale source version
    $ws.Range("P8").Value = $warning

    # Widen the two columns that now hold long URLs / file names.
    $ws.Columns.Item(9).ColumnWidth = 39.1666666666667
    $ws.Columns.Item(16).ColumnWidth = 39.1666666666667
}

Update-LocaleSheet -SheetName "zh-cn" -LocaleRepo "ol-test0-zhcn" `
    -TargetXlf "af351dea-ea3d-4126-985b-09a31b54063f.1d41d7c27e935a8ee1ffbaa08ff96b92b4ae6eb2.zh-cn.xlf" `
    -HandbackDateTime "2016-08-29 22:45:35"

Update-LocaleSheet -SheetName "de-de" -LocaleRepo "ol-test0-dede" `
    -TargetXlf "af351dea-ea3d-4126-985b-09a31b54063f.1d41d7c27e935a8ee1ffbaa08ff96b92b4ae6eb2.de-de.xlf" `
    -HandbackDateTime "2016-08-29 22:45:42"

Write-Output "Generate Report for Handback: done"
